$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 421 is fully populated across columns A-M (no N) and is used purely as a
# style/format donor so the five new rows pick up the same cell styles
# (s="3" for the timestamp column, s="4" for the percentage column, s="2"
# elsewhere) that the rest of the sheet already uses.

function Add-Row {
    param(
        [int]$RowNum,
        [string[]]$Columns,
        [hashtable]$Values
    )

    foreach ($col in $Columns) {
        $ws.Range("$col" + "421").Copy()
        $ws.Range("$col" + "$RowNum").PasteSpecial(-4122)
    }

    foreach ($col in $Columns) {
        $ws.Range("$col" + "$RowNum").Value = $Values[$col]
    }
}

Add-Row 422 @("A","B","C","D","E","F","G","H","I","J","L","M") @{
    A = 44246.689904652776
    B = "PK-Seutu (Helsinki, Espoo, Vantaa)"
    C = "31-35 v"
    D = "Mies"
    E = 11
    F = "Työntekijä / palkollinen"
    G = 1
    H = "Full stack"
    I = "Noin 50/50 hybridimalli"
    J = 7000
    L = "Kyllä"
    M = "Mavericks"
}

Add-Row 423 @("A","B","C","D","E","F","G","H","I","J","K","L","M") @{
    A = 44246.690365104165
    B = "PK-Seutu (Helsinki, Espoo, Vantaa)"
    C = "31-35 v"
    D = "Mies"
    E = 12
    F = "Työntekijä / palkollinen"
    G = 1
    H = "full-stack"
    I = "Pääosin tai kokonaan etätyö"
    J = 8000
    K = 95000
    L = "Kyllä"
    M = "Mavericks"
}

Add-Row 424 @("A","B","C","D","E","F","G","H","I","J","K","L") @{
    A = 44246.69231409722
    B = "Tampere"
    C = "41-45 v"
    D = "Mies"
    E = 22
    F = "Työntekijä / palkollinen"
    G = 0.8
    H = "ohjelmistokehittäjä (backend) / arkkitehti"
    I = "Pääosin tai kokonaan etätyö"
    J = 4700
    K = 58750
    L = "Ei"
}

Add-Row 425 @("A","B","C","D","E","F","G","H","I","J","K","L") @{
    A = 44246.693534756945
    B = "PK-Seutu (Helsinki, Espoo, Vantaa)"
    C = "36-40 v"
    D = "Mies"
    E = 2
    F = "Työntekijä / palkollinen"
    G = 1
    H = "WordPress-kehittäjä"
    I = "Noin 50/50 hybridimalli"
    J = 3000
    K = 37500
    L = "Ei"
}

Add-Row 426 @("A","B","C","D","E","F","G","H","I","J","K","M") @{
    A = 44246.693921655096
    B = "Tampere"
    C = "31-35 v"
    D = "mies"
    E = 5
    F = "Työntekijä / palkollinen"
    G = 1
    H = "Data scientist"
    I = "Pääosin tai kokonaan etätyö"
    J = 4300
    K = 53750
    M = "Wapice"
}
